# Applies the cryptos.xlsx price/volume update described in the commit
# "Updated cryptos list on Sat Feb 24 16:11:22 UTC 2024 with GitHub Actions".
# Column D (Price) holds text-formatted numbers (thousand separators use
# dots, e.g. 51.130.65) so numeric-looking values are written with a leading
# apostrophe to force Excel to keep them as text instead of coercing to a
# real number - this matches the original inlineStr storage in the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "51.130.65"
$ws.Range("E2").Value = "  +0.58%  "

# Row 3
$ws.Range("D3").Value = "2.959.41"
$ws.Range("E3").Value = "  +1.23%  "

# Row 4
$ws.Range("E4").Value = "  -0.15%  "

# Row 5
$ws.Range("D5").Value = "'379.91"
$ws.Range("E5").Value = "  +2.56%  "

# Row 6
$ws.Range("D6").Value = "'102.15"
$ws.Range("E6").Value = "  +1.97%  "

# Row 7
$ws.Range("E7").Value = "  +2.55%  "

# Row 8
$ws.Range("E8").Value = "  -0.08%  "

# Row 9
$ws.Range("D9").Value = "'0.591"
$ws.Range("E9").Value = "  +3.06%  "

# Row 10
$ws.Range("D10").Value = "'36.51"
$ws.Range("E10").Value = "  +2.72%  "

# Row 11
$ws.Range("E11").Value = "  -0.39%  "

# Row 12
$ws.Range("D12").Value = "'0.0858"
$ws.Range("E12").Value = "  +2.88%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'7.81"
$ws.Range("E13").Value = "  +6.46%  "

# Row 14
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "3.419.36"
$ws.Range("E14").Value = "  +1.39%  "

# Row 15
$ws.Range("D15").Value = "'18.28"
$ws.Range("E15").Value = "  +2.77%  "

# Row 16
$ws.Range("D16").Value = "2.952.48"
$ws.Range("E16").Value = "  +0.95%  "

# Row 17
$ws.Range("D17").Value = "'11.22"
$ws.Range("E17").Value = "  -2.28%  "

# Row 18
$ws.Range("D18").Value = "'0.997"
$ws.Range("E18").Value = "  +4.93%  "

# Row 19
$ws.Range("D19").Value = "51.160.69"
$ws.Range("E19").Value = "  +0.81%  "

# Row 20
$ws.Range("D20").Value = "'3.18"
$ws.Range("E20").Value = "  +2.48%  "

# Row 21
$ws.Range("D21").Value = "'12.53"
$ws.Range("E21").Value = "  +2.50%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0960"
$ws.Range("E22").Value = "  +1.78%  "

# Row 23
$ws.Range("D23").Value = "'70.14"
$ws.Range("E23").Value = "  +3.07%  "

# Row 24
$ws.Range("D24").Value = "'266.75"
$ws.Range("E24").Value = "  +1.68%  "

# Row 25
$ws.Range("D25").Value = "'3.20"
$ws.Range("E25").Value = "  +5.41%  "

# Row 26
$ws.Range("D26").Value = "'7.81"
$ws.Range("E26").Value = "  -2.35%  "

# Row 27
$ws.Range("D27").Value = "'7.37"
$ws.Range("E27").Value = "  +1.38%  "

# Row 28
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.05%  "

# Row 29
$ws.Range("D29").Value = "'25.89"
$ws.Range("E29").Value = "  +2.48%  "

# Row 30
$ws.Range("D30").Value = "'0.164"
$ws.Range("E30").Value = "  -0.37%  "

# Row 31
$ws.Range("E31").Value = "  +0.86%  "

# Row 32
$ws.Range("D32").Value = "'10.30"
$ws.Range("E32").Value = "  +4.84%  "

# Row 33
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").Value = "'34.47"
$ws.Range("E33").Value = "  +6.72%  "

# Row 34
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").Value = "'51.22"
$ws.Range("E34").Value = "  +1.69%  "

# Row 35
$ws.Range("D35").Value = "'2.06"
$ws.Range("E35").Value = "  +1.43%  "

# Row 36
$ws.Range("E36").Value = "  -0.49%  "

# Row 37
$ws.Range("E37").Value = "  -0.20%  "

# Row 38
$ws.Range("D38").Value = "'3.27"
$ws.Range("E38").Value = "  +4.94%  "

# Row 39
$ws.Range("E39").Value = "  +2.32%  "

# Row 40
$ws.Range("D40").Value = "'1.83"
$ws.Range("E40").Value = "  +4.02%  "

# Row 41
$ws.Range("D41").Value = "'16.49"
$ws.Range("E41").Value = "  +2.90%  "

# Row 42
$ws.Range("E42").Value = "  +3.42%  "

# Row 43
$ws.Range("D43").Value = "'125.16"
$ws.Range("E43").Value = "  +5.14%  "

# Row 44
$ws.Range("D44").Value = "'3.58"
$ws.Range("E44").Value = "  +11.54%  "

# Row 45
$ws.Range("D45").Value = "'21.42"
$ws.Range("E45").Value = "  +3.56%  "

# Row 46
$ws.Range("E46").Value = "  +4.40%  "

# Row 47
$ws.Range("E47").Value = "  -0.44%  "

# Row 48
$ws.Range("D48").Value = "'0.269"
$ws.Range("E48").Value = "  -0.08%  "

# Row 49
$ws.Range("D49").Value = "2.031.06"
$ws.Range("E49").Value = "  +3.02%  "

# Row 50
$ws.Range("B50").Value = "WOONetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D50").Value = "'0.521"
$ws.Range("E50").Value = "  +14.57%  "

# Row 51
$ws.Range("B51").Value = "BEAM"
$ws.Range("C51").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D51").Value = "'0.0321"
$ws.Range("E51").Value = "  -0.56%  "
